# Intermediate work on excel missing data persist
#
# The sheet gains two new columns: a "Nazev" column at the far left (A)
# and an "Obrazek" column between the old "Katalog c." and "Popisek"
# columns. So the old layout
#   A: Katalog c.   B: Popisek
# becomes
#   A: Nazev   B: Katalog c.   C: Obrazek   D: Popisek
# with C2 left as an empty (but still styled, e.g. for a future picture)
# placeholder cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the old "Popisek" column (B) two slots to the right, into D,
# preserving its values/styles.
$ws.Range("B1").Copy($ws.Range("D1"))
$ws.Range("B2").Copy($ws.Range("D2"))
$ws.Range("B3").Copy($ws.Range("D3"))

# Stake out the new "Obrazek" placeholder cell C2: give it the same
# style as the numeric catalog cell above/left of it, but no content.
$ws.Range("A2").Copy($ws.Range("C2"))
$ws.Range("C2").ClearContents()

# Move the old "Katalog c." column (A) one slot to the right, into B,
# preserving its values/styles.
$ws.Range("A1").Copy($ws.Range("B1"))
$ws.Range("A2").Copy($ws.Range("B2"))
$ws.Range("A3").Copy($ws.Range("B3"))

# Column A is now free for the new "Nazev" column; clear any leftover
# formatting first so it starts from a clean default cell.
$ws.Range("A1:A3").Clear()

# Fill in the two brand new header cells.
$ws.Range("A1").Value = "Nazev"
$ws.Range("C1").Value = "Obrazek"

# Match the author's final selection in the sheet.
[void]$ws.Range("C2").Select()
